$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.921.48'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '3.844.76'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '703.53'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.30'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('D7').Value = '3.841.69'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.524'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.32'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').Value = '4.492.93'
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').Value = '3.860.87'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').Value = '70.954.27'
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.17'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.115'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.43'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.80'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '492.37'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.65%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.72'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('E25').Value = '  -1.22%  '
$ws.Range('E26').Value = '  -3.10%  '
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('E28').Value = '  -3.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.17'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.85%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.40'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('D35').Value = '3.801.04'
$ws.Range('E35').Value = '  +0.94%  '
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  -0.56%  '
$ws.Range('E39').Value = '  +5.70%  '
$ws.Range('E40').Value = '  +6.62%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('E42').Value = '  -5.53%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.000313'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -8.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '163.11'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('B48').Value = 'Arweave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '43.46'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.20%  '
$ws.Range('B49').Value = 'TheGraph'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.299'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.64'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '407.10'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.54%  '
